$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - use same style as existing header (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from H1 (bold, border, centered) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2-7
$values = @{
    2 = @(4, 7)
    3 = @(1, 4)
    4 = @(6, 7)
    5 = @(3, 5)
    6 = @(1, 2)
    7 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
